$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.159832119941711
$ws.Range("B1").Value = 2.755426406860352
$ws.Range("C1").Value = 6.840989589691162
$ws.Range("D1").Value = 1.973016023635864
$ws.Range("E1").Value = 1.066205024719238
